# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Totales / contadores del encabezado ---
$ws.Range("E11").Value = 170820      # VALOR MORA total
$ws.Range("C13").Value = 3           # Cant. Trabajadores
$ws.Range("F13").Value = 1           # Cant. Periodos

# --- La fila 18 pasa a ser la ultima fila de la tabla: toma el formato de
# cierre (borde inferior) que tenia la antigua ultima fila (24) ---
$ws.Range("B24:J24").Copy()
$ws.Range("B18:J18").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Tabla de detalle: actualiza los 3 trabajadores que quedan (periodo 2508) ---
# Fila 16 - NATALIA JUDITH GARCIA NARVAEZ conserva su documento, cambia el periodo
$ws.Range("E16").Value = "2508"

# Fila 17 - ahora BENJAMIN PARRA GUZMAN (antes estaba en la fila 23)
$ws.Range("C17").Value = "10887159"
$ws.Range("D17").Value = "BENJAMIN PARRA GUZMAN"
$ws.Range("E17").Value = "2508"

# Fila 18 - ahora MANUEL SEGUNDO COGOLLO PEREZ (antes estaba en la fila 24)
$ws.Range("C18").Value = "1073988947"
$ws.Range("D18").Value = "MANUEL SEGUNDO COGOLLO PEREZ"
$ws.Range("E18").Value = "2508"

# Elimina las filas de los periodos/trabajadores anteriores que ya no aplican
# (filas 19-24 del estado anterior); las filas de firma se recorren hacia
# arriba automaticamente (antes 29-30, ahora 23-24).
$ws.Rows("19:24").Delete()
